# Update cryptos list values as per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells may contain values that look numeric (e.g. "105.11").
# Format them as Text first so Excel keeps them as literal strings (matching the
# "42.003.72"-style thousand-grouped values already in the sheet), then clear the
# temporary number format so no extra cell styling is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.129.36"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.302.58"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "311.67"
$ws.Range("E5").Value = "  -4.25%  "
$ws.Range("D6").Value = "105.11"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").Value = "40.24"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "0.0914"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "0.973"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "15.51"
$ws.Range("E15").Value = "  -6.28%  "
$ws.Range("D16").Value = "2.650.19"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "2.303.27"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "42.001.99"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").Value = "  -5.37%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "74.43"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  -6.98%  "
$ws.Range("D23").Value = "259.38"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "9.30"
$ws.Range("E25").Value = "  -7.45%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "10.98"
$ws.Range("E27").Value = "  -4.16%  "
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").Value = "22.75"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "35.66"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").Value = "164.25"
$ws.Range("E31").Value = "  -6.30%  "
$ws.Range("D32").Value = "0.0896"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "2.92"
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("E36").Value = "  +11.49%  "
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  -5.31%  "
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("D41").Value = "71.59"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "98.03"
$ws.Range("E42").Value = "  +8.86%  "
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("D47").Value = "112.55"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").Value = "74.04"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.561.53"
$ws.Range("E51").Value = "  -0.11%  "

$dRange.ClearFormats()
